# Refresh the "cryptos" list snapshot (price + 1h volume change)
# for every coin row, plus the row-32/row-33 swap (PancakeSwap <->
# EthereumClassic) exactly as produced by the scraping job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.696.62"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.632.10"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "2.640.56"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.43%  "
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "3.091.33"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "58.692.35"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "2.643.99"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "348.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  -4.07%  "
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.987"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "279.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0983"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "1.988.28"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.10%  "
